$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-38 down to 6-39.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly price-record data.
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44515
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100112027
$ws.Range("G5").Value = "Melón"
$ws.Range("H5").Value = "Tuna"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 26000
$ws.Range("M5").Value = 25500
$ws.Range("N5").Value = "$/caja 18 unidades"
$ws.Range("O5").Value = "Provincia de Copiapó"
$ws.Range("P5").Value = 1417
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"
